$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CRM standard value used for computing "% off" for batch 141 samples
# run on 20210314 (rows 12-19) was entered incorrectly (incrementing by 1
# each row). Correct it to the single actual CRM value for that batch.
$ws.Range("C12:C19").Value = 2234.0700000000002

# Update the sheet view selection to match the saved workbook state
$ws.Range("B18").Select()
